{"js": "// Edit 1: wrap the bold run-pair \"the least amount of \" + \"dimensions\"\n// (inside the PCA-definition paragraph) with a <w:proofErr gramStart/.../gramEnd/>\n// pair, exactly as Word's grammar checker would mark it.\n{\n  const results = context.document.body.search(\"the least amount of dimensions\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find 'the least amount of dimensions' text\");\n  }\n\n  const hit = results.items[0];\n  const paragraph = hit.paragraphs.getFirst();\n  const pRange = paragraph.getRange();\n  pRange.load(\"text\");\n  await context.sync();\n\n  const ooxml1 = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">\n        <w:body>\n          <w:p w14:paraId=\"764A66D1\" w14:textId=\"4681F625\" w:rsidR=\"00F50D5A\" w:rsidRPr=\"00220D68\" w:rsidRDefault=\"00F50D5A\" w:rsidP=\"00F50D5A\">\n            <w:pPr>\n              <w:rPr>\n                <w:rFonts w:cstheme=\"minorHAnsi\"/>\n                <w:b/>\n                <w:bCs/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:lang w:val=\"en-CA\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r w:rsidRPr=\"00220D68\">\n              <w:rPr>\n                <w:rFonts w:cstheme=\"minorHAnsi\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:lang w:val=\"en-CA\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">PCA is used to reduce the dimensionality of a dataset (see overfitting &#8211; lack of generalization), using a transformation that preserves the most variance in the data using </w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramStart\"/>\n            <w:r w:rsidRPr=\"00220D68\">\n              <w:rPr>\n                <w:rFonts w:cstheme=\"minorHAnsi\"/>\n                <w:b/>\n                <w:bCs/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:lang w:val=\"en-CA\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">the least amount of </w:t>\n            </w:r>\n            <w:r w:rsidR=\"00A60877\" w:rsidRPr=\"00220D68\">\n              <w:rPr>\n                <w:rFonts w:cstheme=\"minorHAnsi\"/>\n                <w:b/>\n                <w:bCs/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:lang w:val=\"en-CA\"/>\n              </w:rPr>\n              <w:t>dimensions</w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramEnd\"/>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\n  pRange.insertOoxml(ooxml1, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Edit 2: the empty paragraph right after \"... linear combinations of\n// principal components\" gets filled with \"(in slides, X is the point that\n// has to be flattened)\" \u2014 \"in\" is wrapped in gramStart/gramEnd, and the\n// paragraph-mark formatting is simplified (no more rFonts/shd, just sz/\n// szCs/lang).\n{\n  const results2 = context.document.body.search(\"linear combinations of principal components\", { matchCase: true });\n  results2.load(\"items\");\n  await context.sync();\n\n  if (results2.items.length === 0) {\n    throw new Error(\"Could not find 'linear combinations of principal components' text\");\n  }\n\n  const hit2 = results2.items[0];\n  const paragraph2 = hit2.paragraphs.getFirst();\n  const nextParagraph = paragraph2.getNext();\n  const pRange2 = nextParagraph.getRange();\n  pRange2.load(\"text\");\n  await context.sync();\n\n  const ooxml2 = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">\n        <w:body>\n          <w:p w14:paraId=\"38302BE6\" w14:textId=\"2B10E468\" w:rsidR=\"00184518\" w:rsidRPr=\"00220D68\" w:rsidRDefault=\"00184518\" w:rsidP=\"00F50D5A\">\n            <w:pPr>\n              <w:rPr>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:lang w:val=\"en-CA\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:cstheme=\"minorHAnsi\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n                <w:lang w:val=\"en-CA\"/>\n              </w:rPr>\n              <w:t>(</w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramStart\"/>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:cstheme=\"minorHAnsi\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n                <w:lang w:val=\"en-CA\"/>\n              </w:rPr>\n              <w:t>in</w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramEnd\"/>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:cstheme=\"minorHAnsi\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n                <w:lang w:val=\"en-CA\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\"> slides, </w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:lang w:val=\"en-CA\"/>\n              </w:rPr>\n              <w:t>X is the point that has to be flattened</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:cstheme=\"minorHAnsi\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n                <w:lang w:val=\"en-CA\"/>\n              </w:rPr>\n              <w:t>)</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\n  pRange2.insertOoxml(ooxml2, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Edit 1: wrap the bold run-pair \"the least amount of \" + \"dimensions\"\n# (inside the PCA-definition paragraph) with a <w:proofErr gramStart/.../gramEnd/>\n# pair, exactly as Word's grammar checker would mark it.\n$d = $word.ActiveDocument\n\n# Locate the host paragraph by scanning $d.Paragraphs (Range.Paragraphs.First,\n# derived from a collapsed Find hit, only spans the matched text \u2014 not the\n# whole paragraph \u2014 so we resolve the real paragraph via the document's\n# paragraph collection instead, which gives the correct full Start/End).\n$hostParaIndex1 = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*the least amount of dimensions*\") {\n        $hostParaIndex1 = $i\n        break\n    }\n}\nif ($hostParaIndex1 -eq 0) {\n    throw \"Could not find 'the least amount of dimensions'\"\n}\n$pRange1 = $d.Paragraphs.Item($hostParaIndex1).Range\n\n$ooxml1 = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">\n        <w:body>\n          <w:p w14:paraId=\"764A66D1\" w14:textId=\"4681F625\" w:rsidR=\"00F50D5A\" w:rsidRPr=\"00220D68\" w:rsidRDefault=\"00F50D5A\" w:rsidP=\"00F50D5A\">\n            <w:pPr>\n              <w:rPr>\n                <w:rFonts w:cstheme=\"minorHAnsi\"/>\n                <w:b/>\n                <w:bCs/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:lang w:val=\"en-CA\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r w:rsidRPr=\"00220D68\">\n              <w:rPr>\n                <w:rFonts w:cstheme=\"minorHAnsi\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:lang w:val=\"en-CA\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">PCA is used to reduce the dimensionality of a dataset (see overfitting &#8211; lack of generalization), using a transformation that preserves the most variance in the data using </w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramStart\"/>\n            <w:r w:rsidRPr=\"00220D68\">\n              <w:rPr>\n                <w:rFonts w:cstheme=\"minorHAnsi\"/>\n                <w:b/>\n                <w:bCs/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:lang w:val=\"en-CA\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">the least amount of </w:t>\n            </w:r>\n            <w:r w:rsidR=\"00A60877\" w:rsidRPr=\"00220D68\">\n              <w:rPr>\n                <w:rFonts w:cstheme=\"minorHAnsi\"/>\n                <w:b/>\n                <w:bCs/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:lang w:val=\"en-CA\"/>\n              </w:rPr>\n              <w:t>dimensions</w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramEnd\"/>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$pRange1.InsertXML($ooxml1)\n\n# Edit 2: the empty paragraph right after \"... linear combinations of\n# principal components\" gets filled with \"(in slides, X is the point that\n# has to be flattened)\" \u2014 \"in\" is wrapped in gramStart/gramEnd, and the\n# paragraph-mark formatting is simplified (no more rFonts/shd, just sz/\n# szCs/lang).\n$hostParaIndex2 = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*linear combinations of principal components*\") {\n        $hostParaIndex2 = $i\n        break\n    }\n}\nif ($hostParaIndex2 -eq 0) {\n    throw \"Could not find 'linear combinations of principal components'\"\n}\n\n$pRange2 = $d.Paragraphs.Item($hostParaIndex2 + 1).Range\n\n$ooxml2 = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">\n        <w:body>\n          <w:p w14:paraId=\"38302BE6\" w14:textId=\"2B10E468\" w:rsidR=\"00184518\" w:rsidRPr=\"00220D68\" w:rsidRDefault=\"00184518\" w:rsidP=\"00F50D5A\">\n            <w:pPr>\n              <w:rPr>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:lang w:val=\"en-CA\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:cstheme=\"minorHAnsi\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n                <w:lang w:val=\"en-CA\"/>\n              </w:rPr>\n              <w:t>(</w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramStart\"/>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:cstheme=\"minorHAnsi\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n                <w:lang w:val=\"en-CA\"/>\n              </w:rPr>\n              <w:t>in</w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramEnd\"/>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:cstheme=\"minorHAnsi\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n                <w:lang w:val=\"en-CA\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\"> slides, </w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:lang w:val=\"en-CA\"/>\n              </w:rPr>\n              <w:t>X is the point that has to be flattened</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:cstheme=\"minorHAnsi\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>\n                <w:lang w:val=\"en-CA\"/>\n              </w:rPr>\n              <w:t>)</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$pRange2.InsertXML($ooxml2)\n"}
